$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F11").Value = "dvdsdvdvsdv"
$ws.Range("N11").Value = "sdcsdv"
$ws.Range("J8").Value = "wvevev"
$ws.Range("J16").Value = "dvwevev"

$ws.Range("J16").Select()
